$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-like number format on Price (D) cells before assigning,
# so values such as '482.49' or '68.373.86' are stored as text, not numbers.
$priceCells = @('D2', 'D3', 'D5', 'D6', 'D7', 'D12', 'D14', 'D15', 'D16', 'D18', 'D20', 'D21', 'D23', 'D24', 'D25', 'D27', 'D28', 'D30', 'D32', 'D34', 'D35', 'D36', 'D37', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D49', 'D51')
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = '68.373.86'
$ws.Range("E2").Value = '  +2.16%  '
$ws.Range("D3").Value = '3.899.94'
$ws.Range("E3").Value = '  +1.16%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '482.49'
$ws.Range("E5").Value = '  +1.94%  '
$ws.Range("D6").Value = '145.61'
$ws.Range("E6").Value = '  +0.59%  '
$ws.Range("D7").Value = '0.622'
$ws.Range("E7").Value = '  -1.27%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -2.62%  '
$ws.Range("E10").Value = '  +8.22%  '
$ws.Range("E11").Value = '  +14.24%  '
$ws.Range("D12").Value = '42.79'
$ws.Range("E12").Value = '  -1.74%  '
$ws.Range("E13").Value = '  +2.17%  '
$ws.Range("D14").Value = '4.522.88'
$ws.Range("E14").Value = '  +0.76%  '
$ws.Range("D15").Value = '14.67'
$ws.Range("E15").Value = '  -0.85%  '
$ws.Range("D16").Value = '3.913.26'
$ws.Range("E16").Value = '  +2.43%  '
$ws.Range("E17").Value = '  -0.38%  '
$ws.Range("D18").Value = '19.78'
$ws.Range("E18").Value = '  -1.61%  '
$ws.Range("E19").Value = '  -2.96%  '
$ws.Range("D20").Value = '68.361.72'
$ws.Range("E20").Value = '  +1.77%  '
$ws.Range("D21").Value = '436.78'
$ws.Range("E21").Value = '  +1.30%  '
$ws.Range("E22").Value = '  -0.55%  '
$ws.Range("D23").Value = '3.38'
$ws.Range("E23").Value = '  +1.67%  '
$ws.Range("D24").Value = '87.94'
$ws.Range("E24").Value = '  -0.52%  '
$ws.Range("D25").Value = '11.55'
$ws.Range("E25").Value = '  +16.57%  '
$ws.Range("E26").Value = '  -0.54%  '
$ws.Range("D27").Value = '10.50'
$ws.Range("E27").Value = '  +4.91%  '
$ws.Range("D28").Value = '38.06'
$ws.Range("E28").Value = '  +0.35%  '
$ws.Range("E29").Value = '  +4.75%  '
$ws.Range("D30").Value = '705.29'
$ws.Range("E30").Value = '  -3.48%  '
$ws.Range("E31").Value = '  -2.47%  '
$ws.Range("D32").Value = '13.36'
$ws.Range("E32").Value = '  -3.57%  '
$ws.Range("D34").Value = '0.0₃0927'
$ws.Range("E34").Value = '  +37.89%  '
$ws.Range("D35").Value = '41.81'
$ws.Range("E35").Value = '  -3.45%  '
$ws.Range("D36").Value = '59.65'
$ws.Range("E36").Value = '  +2.65%  '
$ws.Range("D37").Value = '5.68'
$ws.Range("E37").Value = '  +3.97%  '
$ws.Range("E38").Value = '  -6.29%  '
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("D40").Value = '0.0475'
$ws.Range("E40").Value = '  -1.82%  '
$ws.Range("D41").Value = '3.04'
$ws.Range("E41").Value = '  +9.95%  '
$ws.Range("D42").Value = '3.03'
$ws.Range("E42").Value = '  +4.10%  '
$ws.Range("D43").Value = '2.72'
$ws.Range("E43").Value = '  +6.72%  '
$ws.Range("D44").Value = '0.342'
$ws.Range("E44").Value = '  -1.43%  '
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("E46").Value = '  -0.38%  '
$ws.Range("E47").Value = '  -1.42%  '
$ws.Range("E48").Value = '  -0.52%  '
$ws.Range("D49").Value = '146.23'
$ws.Range("E49").Value = '  +2.02%  '
$ws.Range("E50").Value = '  -2.38%  '
$ws.Range("D51").Value = '2.84'
$ws.Range("E51").Value = '  -1.28%  '

# Restore default (General) number format and cell style so no stray
# style/format metadata is left behind on the workbook.
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "General"
    $ws.Range($cell).Style = "Normal"
}
